# Auto-generated PowerShell Excel COM-interop script
# Applies the crypto price/volume list update described in the commit message
# "Updated cryptos list ... with GitHub Actions": refreshed Price (col D) and
# Volume(1h) (col E) figures, plus a 3-row reshuffle of Polkadot/TRON/Chainlink
# (rows 18-20) with their own new Price/Volume values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "67.110.28"
$ws.Cells.Item(2, 5).Value = "  -3.79%  "
$ws.Cells.Item(3, 4).Value = "3.657.77"
$ws.Cells.Item(3, 5).Value = "  -4.63%  "
$ws.Cells.Item(4, 5).Value = "  +0.00%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "594.01"
$ws.Cells.Item(5, 5).Value = "  -3.46%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "164.75"
$ws.Cells.Item(6, 5).Value = "  -7.14%  "
$ws.Cells.Item(7, 4).Value = "3.657.15"
$ws.Cells.Item(7, 5).Value = "  -4.57%  "
$ws.Cells.Item(8, 5).Value = "  +0.02%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.524"
$ws.Cells.Item(9, 5).Value = "  -0.92%  "
$ws.Cells.Item(10, 5).Value = "  -5.87%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "6.16"
$ws.Cells.Item(11, 5).Value = "  -4.65%  "
$ws.Cells.Item(12, 5).Value = "  -5.12%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "37.24"
$ws.Cells.Item(13, 5).Value = "  -6.54%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.0000238"
$ws.Cells.Item(14, 5).Value = "  -6.50%  "
$ws.Cells.Item(15, 4).Value = "4.267.04"
$ws.Cells.Item(15, 5).Value = "  -4.34%  "
$ws.Cells.Item(16, 4).Value = "3.655.14"
$ws.Cells.Item(16, 5).Value = "  -4.42%  "
$ws.Cells.Item(17, 4).Value = "67.164.67"
$ws.Cells.Item(17, 5).Value = "  -3.81%  "
$ws.Cells.Item(18, 2).Value = "Chainlink"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "17.38"
$ws.Cells.Item(18, 5).Value = "  +4.16%  "
$ws.Cells.Item(19, 2).Value = "Polkadot"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.13"
$ws.Cells.Item(19, 5).Value = "  -5.71%  "
$ws.Cells.Item(20, 2).Value = "TRON"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.114"
$ws.Cells.Item(20, 5).Value = "  -3.92%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "490.86"
$ws.Cells.Item(21, 5).Value = "  -3.41%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "9.09"
$ws.Cells.Item(22, 5).Value = "  -5.67%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.714"
$ws.Cells.Item(23, 5).Value = "  -4.14%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "85.30"
$ws.Cells.Item(24, 5).Value = "  -1.30%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.28"
$ws.Cells.Item(25, 5).Value = "  -7.81%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.0000135"
$ws.Cells.Item(26, 5).Value = "  -6.09%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "12.09"
$ws.Cells.Item(27, 5).Value = "  -4.66%  "
$ws.Cells.Item(28, 5).Value = "  -0.55%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "9.91"
$ws.Cells.Item(29, 5).Value = "  -5.55%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "2.91"
$ws.Cells.Item(30, 5).Value = "  -2.33%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "2.35"
$ws.Cells.Item(31, 5).Value = "  -7.13%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "7.62"
$ws.Cells.Item(32, 5).Value = "  -4.36%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "31.36"
$ws.Cells.Item(33, 5).Value = "  -1.79%  "
$ws.Cells.Item(34, 4).Value = "3.795.64"
$ws.Cells.Item(34, 5).Value = "  -4.34%  "
$ws.Cells.Item(35, 5).Value = "  -6.56%  "
$ws.Cells.Item(36, 4).Value = "3.591.91"
$ws.Cells.Item(36, 5).Value = "  -4.55%  "
$ws.Cells.Item(37, 5).Value = "  +0.09%  "
$ws.Cells.Item(38, 5).Value = "  -6.10%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "5.73"
$ws.Cells.Item(39, 5).Value = "  -6.59%  "
$ws.Cells.Item(40, 5).Value = "  -6.87%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.321"
$ws.Cells.Item(41, 5).Value = "  -4.57%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "430.67"
$ws.Cells.Item(42, 5).Value = "  -11.21%  "
$ws.Cells.Item(43, 5).Value = "  -2.34%  "
$ws.Cells.Item(44, 5).Value = "  -7.26%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.74"
$ws.Cells.Item(45, 5).Value = "  -10.04%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "8.32"
$ws.Cells.Item(46, 5).Value = "  -2.87%  "
$ws.Cells.Item(47, 5).Value = "  +0.00%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "40.35"
$ws.Cells.Item(48, 5).Value = "  -8.14%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "142.03"
$ws.Cells.Item(49, 5).Value = "  +1.89%  "
$ws.Cells.Item(50, 4).Value = "2.739.56"
$ws.Cells.Item(50, 5).Value = "  -6.67%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0344"
$ws.Cells.Item(51, 5).Value = "  -4.57%  "
